# Implementación autorización para rol MEC y ADMIN de peticiones de cliente
#
# 1) Row 24: the "delete user" story is reworded to "delete client" (B24),
#    the "Eliminar un cliente por su ID" text (D24) stays as-is.
# 2) Rows 11, 14 and 24 (client list / client details / delete client):
#    ESTADO column (H) flips from "AUTH" (red) to "FLUTTER" (light blue),
#    meaning those client endpoints are now authorized/implemented.
# 3) The active selection moves to H28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24: reword the user story, keep the "Eliminar un cliente por su ID" text ---
$ws.Range("B24").Value = "Como ADMIN, quiero eliminar un cliente"
$ws.Range("D24").Value = "Eliminar un cliente por su ID"

# --- Flip ESTADO (H) from AUTH/red to FLUTTER/light-blue for rows 11, 14, 24 ---
$flutterColor = 15773696  # RGB(0,176,240) packed as BGR OLE color -> fgColor FF00B0F0

foreach ($row in 11, 14, 24) {
    $cell = $ws.Range("H$row")
    $cell.Value = "FLUTTER"
    $cell.Interior.Color = $flutterColor
}

# --- Move the active selection to H28 ---
[void]$ws.Range("H28").Select()
